$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update marking scheme values (concise_ms csv pattern)
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -9.6
$ws.Range("E12").Value = "70.4/140"
